$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 219

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
